$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Genome
$ws.Range("E3").Value = 0.654937947113012
$ws.Range("F3").Value = 0.736001206806953
$ws.Range("G3").Value = 0.6988788482919688
$ws.Range("H3").Value = 0.8359647619519862
$ws.Range("I3").Value = 0.8661543726220495

# Row 4 - Realm
$ws.Range("E4").Value = 0.7496435699182936
$ws.Range("F4").Value = 0.8000792738195573
$ws.Range("G4").Value = 0.7949409567565489
$ws.Range("H4").Value = 0.8817352337062997
$ws.Range("I4").Value = 0.9234266751942479

# Row 5 - Kingdom
$ws.Range("E5").Value = 0.7237581828631745
$ws.Range("F5").Value = 0.763971894656299
$ws.Range("G5").Value = 0.751247925014661
$ws.Range("H5").Value = 0.841000503172785
$ws.Range("I5").Value = 0.9038902749570614

# Row 6 - Phylum
$ws.Range("E6").Value = 0.5824484243387079
$ws.Range("F6").Value = 0.5225821788052979
$ws.Range("G6").Value = 0.4434677292712882
$ws.Range("H6").Value = 0.689064135690569
$ws.Range("I6").Value = 0.8298999850734716

# Row 7 - Class
$ws.Range("C7").Value = 34
$ws.Range("D7").Value = 5845
$ws.Range("E7").Value = 0.478005496950284
$ws.Range("F7").Value = 0.4561834249019896
$ws.Range("G7").Value = 0.3802516783727736
$ws.Range("H7").Value = 0.5895970154962428
$ws.Range("I7").Value = 0.7962965043201163

# Row 8 - Order
$ws.Range("C8").Value = 48
$ws.Range("D8").Value = 5838
$ws.Range("E8").Value = 0.4435437848486384
$ws.Range("F8").Value = 0.479821233466352
$ws.Range("G8").Value = 0.3832090333580366
$ws.Range("H8").Value = 0.5461513139618244
$ws.Range("I8").Value = 0.788353626543978

# Row 9 - Family
$ws.Range("C9").Value = 102
$ws.Range("D9").Value = 5990
$ws.Range("E9").Value = 0.3042229653832192
$ws.Range("F9").Value = 0.3516626362800118
$ws.Range("G9").Value = 0.1680898355702084
$ws.Range("H9").Value = 0.3428587124705174
$ws.Range("I9").Value = 0.732309402087617

# Row 10 - Genus
$ws.Range("C10").Value = 360
$ws.Range("D10").Value = 4673
$ws.Range("E10").Value = 0.360007795286641
$ws.Range("F10").Value = 0.2956360122035187
$ws.Range("G10").Value = 0.06818398303961166
$ws.Range("H10").Value = 0.06206717053352025
$ws.Range("I10").Value = 0.6560871281898071
